$d = $word.ActiveDocument

# 1. The run that hosts the image drawing gets an explicit <w:rPr><w:noProof/></w:rPr>
#    (Word sets this automatically on pictures so the proofing tools skip them).
$shape = $d.InlineShapes.Item(1)
$shape.Range.Font.NoProofing = 1

# 2. Append the new "Aula 05" content after the last paragraph in the body.
#    We build the exact WordprocessingML for the new paragraphs (including the
#    blank separator paragraphs and the spell-check proofErr markers Word
#    inserts around the English terms) and inject it via Range.InsertXML so the
#    resulting markup matches precisely what Word itself would produce.
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xmlFragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r><w:t>O char tem a característica variável de acordo com a tabela ASCII</w:t></w:r>
  <w:r><w:t>, sempre é delimitada com uma aspa só, EXEMPLO: ‘a’</w:t></w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r><w:t xml:space="preserve">Quando for realizar um </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>println</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> e utilizar “” no começo da operação ele considera como uma </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>String</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>.</w:t></w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r><w:t xml:space="preserve">A </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>String</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> é um objeto, por isso não é com letra maiúscula, classes primitivas são todas com letra minúscula.</w:t></w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$endRng.InsertXML($xmlFragment)
